# Add North Puyallup River
#
# Inserts a new row for "Westside Road to North Puyallup" into the
# "Hike Difficulties" table, just before the previously-last row
# ("White River to Sunrise"), which shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item(1)

# Grow the table by one row. The runtime appends the new (blank) row
# at the bottom of the table range.
$newListRow = $tbl.ListRows.Add()

$lastRow = $tbl.Range.Rows.Count            # header + data rows
$newRowIndex = $lastRow                     # row number of the newly appended (blank) row
$insertRowIndex = $newRowIndex - 1          # row number where the new data should actually live

# Move the previous last data row (currently sitting just above the
# new blank row) down into the newly appended row, to make room for
# the new entry at $insertRowIndex.
$ws.Range("A" + $insertRowIndex + ":D" + $insertRowIndex).Copy()
$ws.Range("A" + $newRowIndex).PasteSpecial()
$excel.CutCopyMode = $false

# Fill in the new hike data at the freed-up row.
$ws.Range("A" + $insertRowIndex).Value = "Westside Road to North Puyallup"
$ws.Range("B" + $insertRowIndex).Formula = "=22.5/2"
$ws.Range("C" + $insertRowIndex).Value = 3800
$ws.Range("D" + $insertRowIndex).Value = "strenuous"

# Restore cursor/selection similar to the edited workbook.
$ws.Range("G46").Select()
